$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (font/fill/border/alignment) of column B (rows 1-30)
# into the new column C, then fill in the "Membro" (member) values.
$ws.Range("B1:B30").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("C1").Value = "Membro"
$ws.Range("C2").Value = "José Castro"
$ws.Range("C3").Value = "José Castro"
$ws.Range("C4").Value = "José Castro"
$ws.Range("C5").Value = "José Castro"
$ws.Range("C6").Value = "José Castro"
$ws.Range("C7").Value = "José Castro"
$ws.Range("C8").Value = "Pedro Domingos"
$ws.Range("C9").Value = "João Luís"
$ws.Range("C10").Value = "João Luís"
$ws.Range("C11").Value = "João Luís"
$ws.Range("C12").Value = "Pedro Domingos"
$ws.Range("C13").Value = "Pedro Domingos"
$ws.Range("C14").Value = "João Luís"
$ws.Range("C15").Value = "José Castro"
$ws.Range("C16").Value = "Pedro Domingos"
$ws.Range("C17").Value = "Pedro Domingos"
$ws.Range("C18").Value = "Todos"
$ws.Range("C19").Value = "Pedro Domingos"
$ws.Range("C20").Value = "José Castro"
$ws.Range("C21").Value = "João Luís"
$ws.Range("C22").Value = "João Luís"
$ws.Range("C23").Value = "João Luís"
$ws.Range("C24").Value = "Pedro Domingos"
$ws.Range("C25").Value = "Pedro Domingos"
$ws.Range("C26").Value = "João Luís"
$ws.Range("C27").Value = "João Luís"
$ws.Range("C28").Value = "Pedro Domingos"
$ws.Range("C29").Value = "Pedro Domingos"
$ws.Range("C30").Value = "José Castro"

# Resize columns A-C to fit their (now wider) content, same as Excel's
# "AutoFit Column Width" after the edit.
$ws.Range("A1:C30").EntireColumn.AutoFit()

# Select the same range the author left selected when saving.
$ws.Range("D3:E7").Select()
